$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new RCI column header
$ws.Range("F1").Value = "RCI"

# Add RCI values for each row
$ws.Range("F2").Value = 0.9
$ws.Range("F3").Value = -1.83
$ws.Range("F4").Value = 2.4
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0.97
$ws.Range("F7").Value = 0.32
$ws.Range("F8").Value = -0.45
$ws.Range("F9").Value = -0.9
$ws.Range("F10").Value = -1.01
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = -1.2
$ws.Range("F13").Value = -0.31

# Set column F width to match diff (target stored width 15.1640625;
# engine quantizes ColumnWidth->stored width to 1/6 steps, so 14.3 lands
# on the closest representable stored width, 15.1666...)
$ws.Columns.Item(6).ColumnWidth = 14.3

# Update selection to match the diff (F14)
$ws.Range("F14").Select()
